$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q11").Value = 414
$ws.Range("R11").Value = 45.967
$ws.Range("S11").Value = 2043
$ws.Range("T11").Value = 226.839

$ws.Range("Q18").Value = 492
$ws.Range("R18").Value = 42.452
$ws.Range("S18").Value = 2361
$ws.Range("T18").Value = 203.717
$ws.Range("U18").Value = 3647.2
$ws.Range("V18").Value = 314.695
$ws.Range("W18").Value = 1135.134
$ws.Range("X18").Value = 97.944

$ws.Range("Q27:T27").ClearContents()

$ws.Range("Q32").Value = 774
$ws.Range("R32").Value = 20.508
$ws.Range("S32").Value = 4180
$ws.Range("T32").Value = 110.751

$ws.Range("S43").Value = 2664
$ws.Range("T43").Value = 648.922
$ws.Range("W43").Value = 1496.58
$ws.Range("X43").Value = 364.551

$ws.Range("U45").Value = 14
$ws.Range("V45").Value = 15.984
$ws.Range("W45").Value = 71
$ws.Range("X45").Value = 81.06

$ws.Range("Q46").Value = 741
$ws.Range("R46").Value = 69.194
$ws.Range("S46").Value = 4973
$ws.Range("T46").Value = 464.377
$ws.Range("U46").Value = 1382.641
$ws.Range("V46").Value = 129.11
$ws.Range("W46").Value = 8195.290000000001
$ws.Range("X46").Value = 765.273

$ws.Range("Q48").Value = 102
$ws.Range("R48").Value = 17.61
$ws.Range("S48").Value = 817
$ws.Range("T48").Value = 141.052
$ws.Range("W48").Value = 1049.485
$ws.Range("X48").Value = 181.189

$ws.Range("Q57").Value = 40
$ws.Range("R57").Value = 30.154
$ws.Range("S57").Value = 435
$ws.Range("T57").Value = 327.921
$ws.Range("U57").Value = 29.038
$ws.Range("V57").Value = 21.89
$ws.Range("W57").Value = 274.356
$ws.Range("X57").Value = 206.821

$ws.Range("Q62").Value = 2650
$ws.Range("R62").Value = 40.598
$ws.Range("S62").Value = 24620
$ws.Range("T62").Value = 377.182
$ws.Range("U62").Value = 1088.981
$ws.Range("V62").Value = 16.683
$ws.Range("W62").Value = 7514.751
$ws.Range("X62").Value = 115.127

$ws.Range("Q66").Value = 5546
$ws.Range("R66").Value = 66.194

$ws.Range("U68").Value = 168.136
$ws.Range("V68").Value = 16.131

$ws.Range("S76").Value = 6072
$ws.Range("T76").Value = 628.549

$ws.Range("Q77").Value = 1
$ws.Range("R77").Value = 2.93
$ws.Range("S77").Value = 25
$ws.Range("T77").Value = 73.26000000000001

$ws.Range("Q83").Value = 27
$ws.Range("R83").Value = 5.468
$ws.Range("S83").Value = 321
$ws.Range("T83").Value = 65.009
$ws.Range("U83").Value = 16.109
$ws.Range("V83").Value = 3.262
$ws.Range("W83").Value = 121.828
$ws.Range("X83").Value = 24.673

$ws.Range("Q85").Value = 2580
$ws.Range("R85").Value = 42.672
$ws.Range("S85").Value = 26151
$ws.Range("T85").Value = 432.521
$ws.Range("W85").Value = 3630.141
$ws.Range("X85").Value = 60.04

$ws.Range("S95").Value = 1010
$ws.Range("T95").Value = 535.468
$ws.Range("W95").Value = 1487.374
$ws.Range("X95").Value = 788.5549999999999

$ws.Range("S101").Value = 2510
$ws.Range("T101").Value = 922.018
$ws.Range("U101").Value = 120.81
$ws.Range("V101").Value = 44.378
$ws.Range("W101").Value = 597.228
$ws.Range("X101").Value = 219.384

$ws.Range("Q102").Value = 32
$ws.Range("R102").Value = 51.12
$ws.Range("S102").Value = 150
$ws.Range("T102").Value = 239.626

$ws.Range("W108").Value = 21.47
$ws.Range("X108").Value = 48.626

$ws.Range("Q122").Value = 631
$ws.Range("R122").Value = 36.825
$ws.Range("S122").Value = 1797
$ws.Range("T122").Value = 104.874
$ws.Range("U122").Value = 276.622
$ws.Range("V122").Value = 16.144
$ws.Range("W122").Value = 1779.702
$ws.Range("X122").Value = 103.864

$ws.Range("S128").Value = 129
$ws.Range("T128").Value = 23.795
$ws.Range("U128").Value = 17.297
$ws.Range("V128").Value = 3.191
$ws.Range("W128").Value = 81.39700000000001
$ws.Range("X128").Value = 15.014

$ws.Range("S137").Value = 16686
$ws.Range("T137").Value = 440.885

$ws.Range("Q138").Value = 504
$ws.Range("R138").Value = 49.428
$ws.Range("S138").Value = 2870
$ws.Range("T138").Value = 281.463
$ws.Range("W138").Value = 263.932
$ws.Range("X138").Value = 25.884

$ws.Range("Q140").Value = 1210
$ws.Range("R140").Value = 62.897
$ws.Range("W140").Value = 6618.185
$ws.Range("X140").Value = 344.022

$ws.Range("S155").Value = 2426
$ws.Range("T155").Value = 444.351

$ws.Range("Q156").Value = 206
$ws.Range("R156").Value = 99.089
$ws.Range("S156").Value = 1179
$ws.Range("T156").Value = 567.1180000000001
$ws.Range("U156").Value = 113.892
$ws.Range("V156").Value = 54.784
$ws.Range("W156").Value = 688.346
$ws.Range("X156").Value = 331.105

$ws.Range("Q162").Value = 1907
$ws.Range("R162").Value = 40.787
$ws.Range("S162").Value = 10744
$ws.Range("T162").Value = 229.795
$ws.Range("U162").Value = 39.845
$ws.Range("V162").Value = 0.852
$ws.Range("W162").Value = 990.14
$ws.Range("X162").Value = 21.177

$ws.Range("Q166").Value = 323
$ws.Range("R166").Value = 31.983
$ws.Range("U166").Value = 175.722
$ws.Range("V166").Value = 17.399

$ws.Range("Q181").Value = 1751
$ws.Range("R181").Value = 25.793
$ws.Range("S181").Value = 22520
$ws.Range("T181").Value = 331.733
$ws.Range("W181").Value = 13757.061
$ws.Range("X181").Value = 202.649

$ws.Range("Q182").Value = 23069
$ws.Range("R182").Value = 69.694
$ws.Range("S182").Value = 125220
$ws.Range("T182").Value = 378.305
$ws.Range("U182").Value = 1134
$ws.Range("V182").Value = 3.426
$ws.Range("W182").Value = 23403
$ws.Range("X182").Value = 70.703
